$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue 'D2' '60.957.53'
Set-TextValue 'E2' '  -0.65%  '
Set-TextValue 'D3' '3.365.02'
Set-TextValue 'E3' '  -1.15%  '
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '405.31'
Set-TextValue 'E5' '  -1.39%  '
Set-TextValue 'D6' '135.34'
Set-TextValue 'E6' '  +11.35%  '
Set-TextValue 'D7' '0.590'
Set-TextValue 'E7' '  +2.21%  '
Set-TextValue 'E8' '  +0.05%  '
Set-TextValue 'D9' '0.673'
Set-TextValue 'E9' '  +5.79%  '
Set-TextValue 'D10' '0.120'
Set-TextValue 'E10' '  -1.88%  '
Set-TextValue 'D11' '42.56'
Set-TextValue 'E11' '  +3.93%  '
Set-TextValue 'E12' '  -0.87%  '
Set-TextValue 'D13' '3.884.04'
Set-TextValue 'E13' '  -1.37%  '
Set-TextValue 'D14' '8.31'
Set-TextValue 'E14' '  -0.61%  '
Set-TextValue 'D15' '19.53'
Set-TextValue 'E15' '  +0.41%  '
Set-TextValue 'D16' '3.386.75'
Set-TextValue 'E16' '  +0.03%  '
Set-TextValue 'D17' '60.892.88'
Set-TextValue 'E17' '  -0.62%  '
Set-TextValue 'E18' '  -0.12%  '
Set-TextValue 'D19' '10.94'
Set-TextValue 'E19' '  +1.17%  '
Set-TextValue 'D20' '0.0000127'
Set-TextValue 'E20' '  +4.18%  '
Set-TextValue 'E21' '  -3.22%  '
Set-TextValue 'D22' '83.57'
Set-TextValue 'E22' '  +9.53%  '
Set-TextValue 'D23' '312.28'
Set-TextValue 'E23' '  +4.86%  '
Set-TextValue 'D24' '12.68'
Set-TextValue 'E24' '  -0.84%  '
Set-TextValue 'D25' '3.12'
Set-TextValue 'E25' '  -0.37%  '
Set-TextValue 'D26' '4.76'
Set-TextValue 'E26' '  +11.61%  '
Set-TextValue 'D27' '8.33'
Set-TextValue 'E27' '  +9.88%  '
Set-TextValue 'D28' '29.40'
Set-TextValue 'E28' '  -3.92%  '
Set-TextValue 'D29' '7.43'
Set-TextValue 'E29' '  -7.48%  '
Set-TextValue 'E30' '  +0.63%  '
Set-TextValue 'E31' '  +0.74%  '
Set-TextValue 'E32' '  +0.03%  '
Set-TextValue 'D33' '11.27'
Set-TextValue 'E33' '  -0.67%  '
Set-TextValue 'D34' '41.13'
Set-TextValue 'E34' '  -2.72%  '
Set-TextValue 'E35' '  -1.49%  '
Set-TextValue 'D36' '0.0479'
Set-TextValue 'E36' '  +0.58%  '
Set-TextValue 'D37' '51.83'
Set-TextValue 'E37' '  -0.98%  '
Set-TextValue 'E38' '  -0.01%  '
Set-TextValue 'E39' '  -3.03%  '
Set-TextValue 'D40' '2.91'
Set-TextValue 'E40' '  -3.10%  '
Set-TextValue 'D41' '137.23'
Set-TextValue 'E41' '  +2.70%  '
Set-TextValue 'D42' '1.98'
Set-TextValue 'E42' '  +0.93%  '
Set-TextValue 'E43' '  +0.65%  '
Set-TextValue 'B44' 'TheGraph'
Set-TextValue 'C44' 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue 'D44' '0.290'
Set-TextValue 'E44' '  +3.26%  '
Set-TextValue 'B45' 'NEARProtocol'
Set-TextValue 'C45' 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D45' '4.02'
Set-TextValue 'E45' '  +3.27%  '
Set-TextValue 'D46' '16.61'
Set-TextValue 'E46' '  -2.69%  '
Set-TextValue 'E47' '  +1.67%  '
Set-TextValue 'D48' '21.32'
Set-TextValue 'E48' '  -1.42%  '
Set-TextValue 'D49' '2.117.47'
Set-TextValue 'E49' '  -3.79%  '
Set-TextValue 'D50' '2.30'
Set-TextValue 'E50' '  -4.25%  '
Set-TextValue 'E51' '  +0.01%  '
